$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 4.204118999999999
$ws.Range("H2").Value = 12.612357
$ws.Range("I2").Value = 0.01983154129720676
$ws.Range("J2").Value = 0.01983154129720676
$ws.Range("M2").Value = 0.74396
$ws.Range("N2").Value = 2.23188
$ws.Range("O2").Value = 0.006259003216804254
$ws.Range("P2").Value = 0.006259003216804255
$ws.Range("Q2").Value = 3.127696371239999
$ws.Range("R2").Value = 28.14926734116
$ws.Range("S2").Value = 0.0001241256807734035
$ws.Range("T2").Value = 0.0001241256807734035
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 4.204118999999999
$ws.Range("H3").Value = 12.612357
$ws.Range("I3").Value = 0.01983154129720676
$ws.Range("J3").Value = 0.01983154129720676
$ws.Range("M3").Value = 88.14978533333333
$ws.Range("O3").Value = 0.7416121699579786
$ws.Range("P3").Value = 0.7416121699579786
$ws.Range("Q3").Value = 370.5921873657879
$ws.Range("R3").Value = 3335.329686292092
$ws.Range("S3").Value = 0.01470731237503277
$ws.Range("T3").Value = 0.01470731237503277
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 4.204118999999999
$ws.Range("H4").Value = 12.612357
$ws.Range("I4").Value = 0.01983154129720676
$ws.Range("J4").Value = 0.01983154129720676
$ws.Range("M4").Value = 29.76859933333333
$ws.Range("N4").Value = 89.305798
$ws.Range("O4").Value = 0.2504459365921425
$ws.Range("P4").Value = 0.2504459365921425
$ws.Range("Q4").Value = 125.150734060654
$ws.Range("R4").Value = 1126.356606545886
$ws.Range("S4").Value = 0.004966728934244698
$ws.Range("T4").Value = 0.004966728934244699
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 4.204118999999999
$ws.Range("H5").Value = 12.612357
$ws.Range("I5").Value = 0.01983154129720676
$ws.Range("J5").Value = 0.01983154129720676
$ws.Range("M5").Value = 0.2000323333333334
$ws.Range("N5").Value = 0.6000970000000001
$ws.Range("O5").Value = 0.00168289023307462
$ws.Range("P5").Value = 0.00168289023307462
$ws.Range("Q5").Value = 0.840959733181
$ws.Range("R5").Value = 7.568637598629001
$ws.Range("S5").Value = [double]"3.337430715588524E-05"
$ws.Range("T5").Value = [double]"3.337430715588524E-05"
$ws.Range("I6").Value = 0.8539093107807857
$ws.Range("J6").Value = 0.8539093107807858
$ws.Range("M6").Value = 0.74396
$ws.Range("N6").Value = 2.23188
$ws.Range("O6").Value = 0.006259003216804254
$ws.Range("P6").Value = 0.006259003216804255
$ws.Range("Q6").Value = 134.6727928339733
$ws.Range("R6").Value = 1212.05513550576
$ws.Range("S6").Value = 0.005344621123036041
$ws.Range("T6").Value = 0.005344621123036043
$ws.Range("I7").Value = 0.8539093107807857
$ws.Range("J7").Value = 0.8539093107807858
$ws.Range("M7").Value = 88.14978533333333
$ws.Range("O7").Value = 0.7416121699579786
$ws.Range("P7").Value = 0.7416121699579786
$ws.Range("S7").Value = 0.6332695369154604
$ws.Range("T7").Value = 0.6332695369154605
$ws.Range("I8").Value = 0.8539093107807857
$ws.Range("J8").Value = 0.8539093107807858
$ws.Range("M8").Value = 29.76859933333333
$ws.Range("N8").Value = 89.305798
$ws.Range("O8").Value = 0.2504459365921425
$ws.Range("P8").Value = 0.2504459365921425
$ws.Range("Q8").Value = 5388.7580124947
$ws.Range("R8").Value = 48498.8221124523
$ws.Range("S8").Value = 0.2138581171032447
$ws.Range("T8").Value = 0.2138581171032448
$ws.Range("I9").Value = 0.8539093107807857
$ws.Range("J9").Value = 0.8539093107807858
$ws.Range("M9").Value = 0.2000323333333334
$ws.Range("N9").Value = 0.6000970000000001
$ws.Range("O9").Value = 0.00168289023307462
$ws.Range("P9").Value = 0.00168289023307462
$ws.Range("Q9").Value = 36.21016316347156
$ws.Range("R9").Value = 325.8914684712441
$ws.Range("S9").Value = 0.001437035639044465
$ws.Range("T9").Value = 0.001437035639044465
$ws.Range("G10").Value = 26.057747
$ws.Range("H10").Value = 78.173241
$ws.Range("I10").Value = 0.1229188055196976
$ws.Range("J10").Value = 0.1229188055196976
$ws.Range("M10").Value = 0.74396
$ws.Range("N10").Value = 2.23188
$ws.Range("O10").Value = 0.006259003216804254
$ws.Range("P10").Value = 0.006259003216804255
$ws.Range("Q10").Value = 19.38592145812
$ws.Range("R10").Value = 174.47329312308
$ws.Range("S10").Value = 0.0007693491991535239
$ws.Range("T10").Value = 0.000769349199153524
$ws.Range("G11").Value = 26.057747
$ws.Range("H11").Value = 78.173241
$ws.Range("I11").Value = 0.1229188055196976
$ws.Range("J11").Value = 0.1229188055196976
$ws.Range("M11").Value = 88.14978533333333
$ws.Range("O11").Value = 0.7416121699579786
$ws.Range("P11").Value = 0.7416121699579786
$ws.Range("Q11").Value = 2296.984804320311
$ws.Range("R11").Value = 20672.86323888279
$ws.Range("S11").Value = 0.0911580820901057
$ws.Range("T11").Value = 0.0911580820901057
$ws.Range("G12").Value = 26.057747
$ws.Range("H12").Value = 78.173241
$ws.Range("I12").Value = 0.1229188055196976
$ws.Range("J12").Value = 0.1229188055196976
$ws.Range("M12").Value = 29.76859933333333
$ws.Range("N12").Value = 89.305798
$ws.Range("O12").Value = 0.2504459365921425
$ws.Range("P12").Value = 0.2504459365921425
$ws.Range("Q12").Value = 775.7026299723686
$ws.Range("R12").Value = 6981.323669751318
$ws.Range("S12").Value = 0.03078451537316808
$ws.Range("T12").Value = 0.03078451537316808
$ws.Range("G13").Value = 26.057747
$ws.Range("H13").Value = 78.173241
$ws.Range("I13").Value = 0.1229188055196976
$ws.Range("J13").Value = 0.1229188055196976
$ws.Range("M13").Value = 0.2000323333333334
$ws.Range("N13").Value = 0.6000970000000001
$ws.Range("O13").Value = 0.00168289023307462
$ws.Range("P13").Value = 0.00168289023307462
$ws.Range("Q13").Value = 5.212391933819668
$ws.Range("R13").Value = 46.91152740437701
$ws.Range("S13").Value = 0.0002068588572702978
$ws.Range("T13").Value = 0.0002068588572702979
$ws.Range("G14").Value = 0.7081243333333332
$ws.Range("H14").Value = 2.124373
$ws.Range("I14").Value = 0.003340342402309973
$ws.Range("J14").Value = 0.003340342402309974
$ws.Range("M14").Value = 0.74396
$ws.Range("N14").Value = 2.23188
$ws.Range("O14").Value = 0.006259003216804254
$ws.Range("P14").Value = 0.006259003216804255
$ws.Range("Q14").Value = 0.5268161790266666
$ws.Range("R14").Value = 4.741345611239999
$ws.Range("S14").Value = [double]"2.090721384128577E-05"
$ws.Range("T14").Value = [double]"2.090721384128578E-05"
$ws.Range("G15").Value = 0.7081243333333332
$ws.Range("H15").Value = 2.124373
$ws.Range("I15").Value = 0.003340342402309973
$ws.Range("J15").Value = 0.003340342402309974
$ws.Range("M15").Value = 88.14978533333333
$ws.Range("O15").Value = 0.7416121699579786
$ws.Range("P15").Value = 0.7416121699579786
$ws.Range("Q15").Value = 62.4210079726431
$ws.Range("R15").Value = 561.7890717537879
$ws.Range("S15").Value = 0.002477238577379747
$ws.Range("T15").Value = 0.002477238577379747
$ws.Range("G16").Value = 0.7081243333333332
$ws.Range("H16").Value = 2.124373
$ws.Range("I16").Value = 0.003340342402309973
$ws.Range("J16").Value = 0.003340342402309974
$ws.Range("M16").Value = 29.76859933333333
$ws.Range("N16").Value = 89.305798
$ws.Range("O16").Value = 0.2504459365921425
$ws.Range("P16").Value = 0.2504459365921425
$ws.Range("Q16").Value = 21.07986955718377
$ws.Range("R16").Value = 189.718826014654
$ws.Range("S16").Value = 0.0008365751814849684
$ws.Range("T16").Value = 0.0008365751814849686
$ws.Range("G17").Value = 0.7081243333333332
$ws.Range("H17").Value = 2.124373
$ws.Range("I17").Value = 0.003340342402309973
$ws.Range("J17").Value = 0.003340342402309974
$ws.Range("M17").Value = 0.2000323333333334
$ws.Range("N17").Value = 0.6000970000000001
$ws.Range("O17").Value = 0.00168289023307462
$ws.Range("P17").Value = 0.00168289023307462
$ws.Range("Q17").Value = 0.1416477626867778
$ws.Range("R17").Value = 1.274829864181
$ws.Range("S17").Value = [double]"5.621429603972468E-06"
$ws.Range("T17").Value = [double]"5.621429603972469E-06"
